$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Turn the existing "Invalid Key Press Sound" URL plain-text run into a
#    real hyperlink (matches the first hunk of the diff).
# ---------------------------------------------------------------------------
$findRng = $d.Content
$null = $findRng.Find.Execute(
    "http://www.freesound.org/people/Splashdust/sounds/67454/",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $d.Hyperlinks.Add($findRng, "http://www.freesound.org/people/Splashdust/sounds/67454/")

# Locate the paragraph that now holds "Invalid Key Press Sound: <link>" so we
# can anchor the new paragraphs right after it.
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Invalid Key Press Sound:*") {
        $anchorPara = $d.Paragraphs.Item($i)
        break
    }
}

# ---------------------------------------------------------------------------
# Helper: clone the formatting of $anchorPara by inserting a paragraph break
# right after it (this carries over pStyle/numPr/rPr exactly like Word does),
# then type a label run followed by a hyperlink run.
# ---------------------------------------------------------------------------
function Add-SoundParagraph($afterPara, $label, $url) {
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newPara = $afterPara.Next()
    $pr = $newPara.Range
    $pr.Collapse(1)
    $pr.InsertAfter($label)

    $afterLabel = $newPara
    $insPos = $afterLabel.Range.End - 1
    $urlStart = $insPos
    $insRng = $d.Range($insPos, $insPos)
    $insRng.InsertAfter($url)

    $urlEnd = $newPara.Range.End - 1
    $hRng = $d.Range($urlStart, $urlEnd)
    $null = $d.Hyperlinks.Add($hRng, $url)

    return $newPara
}

$p1 = Add-SoundParagraph $anchorPara "Menu selection change sound: " "http://www.freesound.org/people/broumbroum/sounds/50561/"
$p2 = Add-SoundParagraph $p1 "Back To Menu Sound: " "http://www.freesound.org/people/broumbroum/sounds/50557/"
$p3 = Add-SoundParagraph $p2 "Menu Selection Sound: " "http://www.freesound.org/people/broumbroum/sounds/50565/"

# ---------------------------------------------------------------------------
# Final paragraph: plain "Ok" (no hyperlink), cloned the same way.
# ---------------------------------------------------------------------------
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p4 = $p3.Next()
$pr4 = $p4.Range
$pr4.Collapse(1)
$pr4.InsertAfter("Ok")

# ---------------------------------------------------------------------------
# Move the _GoBack bookmark from the old last-run position to the end of the
# new "Ok" paragraph (right before its paragraph mark), matching the diff.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$gbPos = $p4.Range.End - 1
$gbRng = $d.Range($gbPos, $gbPos)
$null = $d.Bookmarks.Add("_GoBack", $gbRng)
